# aggiornamento fino a 20/09/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(375, 44449, 44, 193, 102.1093787199822),
    @(376, 44450, 22, 191, 101.0512504430912),
    @(377, 44451, 42, 207, 109.5162766582193),
    @(378, 44452, 22, 174, 92.05716008951765),
    @(379, 44453, 24, 181, 95.76060905863618),
    @(380, 44454, 16, 188, 99.46405802775472),
    @(381, 44455, 5, 175, 92.58622422796316),
    @(382, 44456, 30, 161, 85.17932628972611),
    @(383, 44457, 16, 155, 82.00494145905309),
    @(384, 44458, 25, 138, 73.01085110547952),
    @(385, 44459, 12, 128, 67.72020972102447)
)

# Source row to copy the date-cell formatting (style) from
$formatSrc = $ws.Range("A374")

foreach ($entry in $data) {
    $r = $entry[0]

    # Replicate column A's date number format / border / font by copying
    # the format of the last existing date cell before writing the value.
    $formatSrc.Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}
